# Update the SRA_data sheet:
#  - Column J ("design_description"): replace "Sequencing performed at None"
#    with "Sequencing performed at Oregon State University Center for
#    Quantitative Life Sciences Genomics Core"
#  - Columns L/M ("filename"/"filename2"): rebuild the read-1/read-2 fastq
#    filenames from the library_ID (column B) instead of the stale
#    "MP_..." placeholder names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SRA_data")

$oldDesc = "Sequencing performed at None"
$newDesc = "Sequencing performed at Oregon State University Center for Quantitative Life Sciences Genomics Core"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $libId = $ws.Cells.Item($r, 2).Value2

    if ([string]::IsNullOrEmpty($libId)) {
        continue
    }

    $designCell = $ws.Cells.Item($r, 10)
    if ($designCell.Value2 -eq $oldDesc) {
        $designCell.Value = $newDesc
    }

    $ws.Cells.Item($r, 12).Value = $libId + "_R1.fastq.gz"
    $ws.Cells.Item($r, 13).Value = $libId + "_R2.fastq.gz"
}
